{"js": "// Rewrites the \"time\" essay into a \"chemistry\" essay:\n// - Title, author, email updated\n// - Body + Summary paragraphs rewritten with new content\n// - A trailing empty paragraph is appended at the very end of the body\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph layout in this document:\n// 0: Title\n// 1: Author name\n// 2: Author email\n// 3: (blank spacer paragraph)\n// 4: Main body paragraph (contains manual line breaks)\n// 5: \"Summary\" heading\n// 6: Summary paragraph\nconst titlePara = paragraphs.items[0];\nconst authorPara = paragraphs.items[1];\nconst emailPara = paragraphs.items[2];\nconst bodyPara = paragraphs.items[4];\nconst summaryPara = paragraphs.items[6];\n\n// Title\ntitlePara.insertText(\n  \"Unraveling the Mysteries of Chemistry: A Journey into the World of Elements and Reactions\",\n  Word.InsertLocation.replace\n);\n\n// Author name\nauthorPara.insertText(\"Dr. Emily Carter\", Word.InsertLocation.replace);\n\n// Author email\nemailPara.insertText(\"emcarter@chemistryeducators.org\", Word.InsertLocation.replace);\n\n// Main body (uses \\u000b for the manual line breaks originally encoded as <w:br/>)\nconst newBodyText =\n  \"Have you ever wondered about the world around you? Why do things change? \" +\n  \"Why are there so many different substances? The answers to these questions lie in \" +\n  \"one of the most intriguing and fundamental subjects--chemistry. \" +\n  \"On our voyage through the captivating realm of chemistry, we will unravel the secrets \" +\n  \"of matter, understand how substances interact, and explore the incredible applications \" +\n  \"of chemistry in our everyday lives.\\u000b\\u000b\" +\n  \"In this extraordinary odyssey, we will explore the basic building blocks of all matter--the elements. \" +\n  \"We will investigate their properties, their bonding behavior, and their arrangements to form different compounds. \" +\n  \"Through demonstrations and hands-on experiments, we will uncover the mysteries of chemical reactions, \" +\n  \"learning how atoms rearrange and energy is transferred. \" +\n  \"The world of chemistry is a place of fascinating phenomena, from the colorful fireworks that light up our skies \" +\n  \"to the complex processes occurring within our bodies.\\u000b\\u000b\" +\n  \"While embarking on this journey of discovery, we will delve into the diverse applications of chemistry. \" +\n  \"From understanding the role of chemistry in fields such as medicine, engineering, and agriculture to comprehending \" +\n  \"the impact of chemistry on our environment, we will appreciate the significance of this science in shaping our world. \" +\n  \"Chemistry provides solutions to real-world problems, enhances our lives, and continues to push the boundaries of human knowledge.\";\nbodyPara.insertText(newBodyText, Word.InsertLocation.replace);\n\n// Summary paragraph\nconst newSummaryText =\n  \"Our exploration of chemistry has unveiled the fundamental concepts and applications of this dynamic science. \" +\n  \"We have journeyed through the world of elements, reactions, and compounds, unraveling the secrets of matter and its transformations. \" +\n  \"Throughout our voyage, we have witnessed the power of chemistry in diverse fields, from medicine to engineering. \" +\n  \"This knowledge equips us with a deeper understanding of the world around us, empowering us to appreciate the intricacies \" +\n  \"of chemical processes and their impact on our lives. \" +\n  \"As we continue our scientific odyssey, we look forward to unraveling even more mysteries of the chemical realm.\";\nsummaryPara.insertText(newSummaryText, Word.InsertLocation.replace);\n\nawait context.sync();\n\n// A new empty paragraph is added at the very end of the document body.\nbody.insertParagraph(\"\", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# Rewrites the \"time\" essay into a \"chemistry\" essay:\n# - Title, author, email updated\n# - Body + Summary paragraphs rewritten with new content\n# - A trailing empty paragraph is appended at the very end of the body\n\n$doc = $word.ActiveDocument\n\n# Paragraph layout in this document (1-based Paragraphs collection):\n# 1: Title\n# 2: Author name\n# 3: Author email\n# 4: (blank spacer paragraph)\n# 5: Main body paragraph (contains manual line breaks)\n# 6: \"Summary\" heading\n# 7: Summary paragraph\n\nfunction Set-ParagraphText($paragraph, $text) {\n    $start = $paragraph.Range.Start\n    $end = $paragraph.Range.End\n    $r = $doc.Range($start, $end)\n    $r.Text = $text\n}\n\n# Title\nSet-ParagraphText $doc.Paragraphs.Item(1) \"Unraveling the Mysteries of Chemistry: A Journey into the World of Elements and Reactions\"\n\n# Author name\nSet-ParagraphText $doc.Paragraphs.Item(2) \"Dr. Emily Carter\"\n\n# Author email\nSet-ParagraphText $doc.Paragraphs.Item(3) \"emcarter@chemistryeducators.org\"\n\n# Main body (uses vertical-tab char for the manual line breaks, i.e. <w:br/>)\n$vt = [char]0x0B\n$newBodyText = \"Have you ever wondered about the world around you? Why do things change? \" +\n  \"Why are there so many different substances? The answers to these questions lie in \" +\n  \"one of the most intriguing and fundamental subjects--chemistry. \" +\n  \"On our voyage through the captivating realm of chemistry, we will unravel the secrets \" +\n  \"of matter, understand how substances interact, and explore the incredible applications \" +\n  \"of chemistry in our everyday lives.$vt$vt\" +\n  \"In this extraordinary odyssey, we will explore the basic building blocks of all matter--the elements. \" +\n  \"We will investigate their properties, their bonding behavior, and their arrangements to form different compounds. \" +\n  \"Through demonstrations and hands-on experiments, we will uncover the mysteries of chemical reactions, \" +\n  \"learning how atoms rearrange and energy is transferred. \" +\n  \"The world of chemistry is a place of fascinating phenomena, from the colorful fireworks that light up our skies \" +\n  \"to the complex processes occurring within our bodies.$vt$vt\" +\n  \"While embarking on this journey of discovery, we will delve into the diverse applications of chemistry. \" +\n  \"From understanding the role of chemistry in fields such as medicine, engineering, and agriculture to comprehending \" +\n  \"the impact of chemistry on our environment, we will appreciate the significance of this science in shaping our world. \" +\n  \"Chemistry provides solutions to real-world problems, enhances our lives, and continues to push the boundaries of human knowledge.\"\nSet-ParagraphText $doc.Paragraphs.Item(5) $newBodyText\n\n# Summary paragraph\n$newSummaryText = \"Our exploration of chemistry has unveiled the fundamental concepts and applications of this dynamic science. \" +\n  \"We have journeyed through the world of elements, reactions, and compounds, unraveling the secrets of matter and its transformations. \" +\n  \"Throughout our voyage, we have witnessed the power of chemistry in diverse fields, from medicine to engineering. \" +\n  \"This knowledge equips us with a deeper understanding of the world around us, empowering us to appreciate the intricacies \" +\n  \"of chemical processes and their impact on our lives. \" +\n  \"As we continue our scientific odyssey, we look forward to unraveling even more mysteries of the chemical realm.\"\nSet-ParagraphText $doc.Paragraphs.Item(7) $newSummaryText\n\n# A new empty paragraph is added at the very end of the document body.\n$last = $doc.Paragraphs.Item($doc.Paragraphs.Count)\n$last.Range.InsertParagraphAfter()\n"}
